# Update "想去人数" (attendance count) values on the 展览 and 全部类型 sheets
# to reflect refreshed data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 3 F 393 -> 401, row 4 F 3059 -> 3081
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 401
$wsExhibit.Range("F4").Value = 3081

# Sheet "全部类型": row 4 F 393 -> 401, row 5 F 3059 -> 3081
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 401
$wsAll.Range("F5").Value = 3081
